$d = $word.ActiveDocument

# --- Paragraph 1 text restructuring -----------------------------------
# Before : "NRG" + "-" + "rkamath"(flagged misspelled) + "-Ending Ma"
# After  : "NRG" + "-rkamath" + "-" + "Ending 2*Ma"   (no spell-flag left
#          on "rkamath", and "Ma" becomes "2*Ma").
#
# Remove the "-rkamath-" span together with its spell-check markers, then
# retype "-rkamath" right after "NRG" and "-" right before "Ending Ma" so
# the run/proofErr bookkeeping gets rebuilt from scratch (mirrors how a
# user re-typing that part of the line would leave the document).

$d.Content.Find.Execute("-rkamath-", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

$d.Content.Find.Execute("NRG", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NRG-rkamath", 2)

$d.Content.Find.Execute("Ending Ma", $true, $false, $false, $false, $false,
                         $true, 1, $false, "-Ending 2*Ma", 2)
